# Generate Report for Archive
#
# The localization status for the two handed-off files moved from
# "Ready for handoff" to "In Translation". That status string is shared
# across the Overview sheet (columns E/F) and the per-locale sheets
# zh-cn / de-de (column C), so a single text replacement applied
# workbook-wide updates every occurrence.
#
# Because the status text got shorter, the Status/zh-cn/de-de columns
# were re-sized (narrower) when the report was regenerated.

$wb = $excel.ActiveWorkbook

# --- 1. Update the status text everywhere it appears in the workbook ---
foreach ($ws in $wb.Worksheets) {
    $ws.Cells.Replace("Ready for handoff", "In Translation")
}

# --- 2. Shrink the now-narrower status columns to match the regenerated report ---
$overview = $wb.Worksheets.Item("Overview")
$overview.Columns.Item(5).ColumnWidth = 12.5   # column E
$overview.Columns.Item(6).ColumnWidth = 12.5   # column F

$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Columns.Item(3).ColumnWidth = 12.5        # column C (Status)

$dede = $wb.Worksheets.Item("de-de")
$dede.Columns.Item(3).ColumnWidth = 12.5        # column C (Status)
